# Update the third row of the flash-card test data:
#   - column B ("amy") becomes the Chinese name "爱丽丝·宝琳"
#   - column A ("amy") becomes "Alice pauline", styled with a custom font
#     (JetBrains Mono, ~9.8pt, RGB 6A8759) like the rest of the new entries
# Set B3 first so the new shared strings are appended in the same order
# as the target workbook (B3's string before A3's string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "爱丽丝·宝琳"
$ws.Range("A3").Value = "Alice pauline"

$ws.Range("A3").Font.Name = "JetBrains Mono"
$ws.Range("A3").Font.Size = 9.8
$ws.Range("A3").Font.Color = 5867370

# Leave the selection on B4, matching the saved view state.
$ws.Range("B4").Select()
